$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column D (shifts the old "Platform" column D -> E) ---
$ws.Columns("D").Insert()

# --- Header row (row 1) ---
$ws.Range("C1").Value = "Problem(s) Solved"
$ws.Range("D1").Value = "Problem(s) Attempted"

# --- New rows at the bottom (6th/7th March entries) ---
$ws.Range("B20").Value = 6
$ws.Range("D20").Value = "Elimination Game (Recursion)"
$ws.Range("E20").Value = "Bosscoder Academy"

$ws.Range("D21").Value = "Elimination Game (Recursion)"
$ws.Range("E21").Value = "LeetCode"

$ws.Range("B22").Value = 7
$ws.Range("C22").Value = "Valid Sudoku"
$ws.Range("E22").Value = "LeetCode"

$ws.Range("C23").Value = "Spiral Matrix"
$ws.Range("E23").Value = "LeetCode"

# --- Column widths (D grew to fit "Elimination Game (Recursion)", E keeps ~old D width) ---
$ws.Range("D1").EntireColumn.ColumnWidth = 34.25
$ws.Range("E1").EntireColumn.ColumnWidth = 16.42

# --- View state: scrolled down, new selection on the newly entered cell ---
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("E24").Select()
